$wb = $excel.ActiveWorkbook

# Sheet1: "Weekly Quantity" - delete rows 11-14
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("A11:B14").EntireRow.Delete() | Out-Null

# Sheet2: "Monthly Trend" - change B7 to 10, delete row 8
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B7").Value = 10
$ws2.Range("A8:B8").EntireRow.Delete() | Out-Null
